$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append row 20 with the new forecast vector values, reusing the
# date-column formatting already applied to the rows above it.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -2.451276118722334
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -1.596682557877005
